# Update Gauteng Covid 19 Data workbook with data for 05/05/2020 (column V)

$wb = $excel.ActiveWorkbook

$district = $wb.Worksheets.Item("District")
$subDistrict = $wb.Worksheets.Item("sub-District")

# ---- District sheet: column V (05/05/2020) ----
# Rows 2-7: plain daily counts (new cells, no existing formatting)
$district.Range("V2").Value = 971
$district.Range("V3").Value = 366
$district.Range("V4").Value = 224
$district.Range("V5").Value = 21
$district.Range("V6").Value = 51
$district.Range("V7").Value = 64

# Rows 8-10: Total Cases / Deaths / Recoveries
$district.Range("V8").Value = 1697
$district.Range("V9").Value = 15
$district.Range("V10").Value = 1036

# Row 11: Active Cases = Total Cases - Deaths - Recoveries, filled across as a
# shared formula (matches how the rest of the row B11:V11 is already built).
$district.Range("B11:V11").Formula = "=B8-B9-B10"

# Rows 12-16: Region recoveries breakdown
$district.Range("V12").Value = 627
$district.Range("V13").Value = 108
$district.Range("V14").Value = 263
$district.Range("V15").Value = 6
$district.Range("V16").Value = 32

# ---- sub-District sheet: column V (05/05/2020) ----
# Rows 2, 11, 19, 28, 33, 38-42 already hold formulas pulling from District!V*,
# so they update automatically once District is recalculated. The remaining
# rows are plain daily counts for each sub-region.
$subDistrict.Range("V3").Value = 133
$subDistrict.Range("V4").Value = 156
$subDistrict.Range("V5").Value = 67
$subDistrict.Range("V6").Value = 75
$subDistrict.Range("V7").Value = 265
$subDistrict.Range("V8").Value = 115
$subDistrict.Range("V9").Value = 67
$subDistrict.Range("V10").Value = 93

$subDistrict.Range("V12").Value = 37
$subDistrict.Range("V13").Value = 31
$subDistrict.Range("V14").Value = 74
$subDistrict.Range("V15").Value = 94
$subDistrict.Range("V16").Value = 83
$subDistrict.Range("V17").Value = 12
$subDistrict.Range("V18").Value = 35

$subDistrict.Range("V20").Value = 32
$subDistrict.Range("V21").Value = 4
$subDistrict.Range("V22").Value = 67
$subDistrict.Range("V23").Value = 56
$subDistrict.Range("V24").Value = 3
$subDistrict.Range("V25").Value = 49
$subDistrict.Range("V26").Value = 7
$subDistrict.Range("V27").Value = 6

$subDistrict.Range("V29").Value = 3
$subDistrict.Range("V30").Value = 14
$subDistrict.Range("V31").Value = 4
$subDistrict.Range("V32").Value = 0

$subDistrict.Range("V34").Value = 36
$subDistrict.Range("V35").Value = 10
$subDistrict.Range("V36").Value = 3
$subDistrict.Range("V37").Value = 2

# ---- View state ----
# District was the active sheet while entering the previous day's figures;
# now that today's column has been added, sub-District becomes the active
# sheet/tab and the selection on each sheet moves to the new column.
$district.Range("V15").Select()
$subDistrict.Select()
$subDistrict.Range("U14").Select()
